$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

$newRow = 24

# Copy formatting (styles + row height) from the row above onto the new row
$ws.Range("A23:E23").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item($newRow).RowHeight = 135

# Populate the new row's values
$ws.Cells.Item($newRow, 1).Value = "IPAIAM0006"
$ws.Cells.Item($newRow, 2).Value = "OPQA-4531||OPQA-4533||OPQA-4547||OPQA-4550"
$ws.Cells.Item($newRow, 3).Value = "Verify that the new STeAM step up authentication modal should include a link to initiate the EndNote supported forgot password flow. || Verify that the `"Sign in to Target Druggability`" modal should match with wireframe || Verify that the target application product overview page should be opened in a separate browser window when User clicks `"Learn more`" in any of the Step up related messages/modals.||Verify that the DRA\IPA application overview page should be opened in a separate browser window when user has a valid session token on the Neon platform"
$ws.Cells.Item($newRow, 4).Value = "Y"
$ws.Cells.Item($newRow, 5).Value = $null

# Update the visible top-left / selected cell, matching the workbook state after the edit
$ws.Application.ActiveWindow.ScrollRow = $newRow
$ws.Range("A$newRow").Select()
